# Applies the FY2022 "dropped COVID revs" refresh:
#  - Rename Table 4.a / Table 4.b to their descriptive names
#  - Tidy a handful of category labels (singular "Tax" instead of "Taxes")
#  - Rename the "Family Services Net Medicaid" category
#  - Refresh a handful of figures that moved after the rerender

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Rename the two CAGR worksheets
# ---------------------------------------------------------------
$wb.Worksheets.Item("Table 4.a").Name = "Revenue CAGR"
$wb.Worksheets.Item("Table 4.b").Name = "Expenditure CAGR"

# ---------------------------------------------------------------
# 2. Category label text fixes (shared across every sheet that
#    references them, so sweep the whole workbook)
# ---------------------------------------------------------------
$labelChanges = @{
    "Individual Income Taxes"        = "Individual Income Tax"
    "Corporate Income Taxes"         = "Corporate Income Tax"
    "Sales Taxes"                    = "Sales Tax"
    "Public Utility Taxes"           = "Public Utility Tax"
    "Family Services Net Medicaid"   = "Healthcare & Fam Ser Net Of Medicaid"
}

foreach ($ws in $wb.Worksheets) {
    foreach ($old in $labelChanges.Keys) {
        [void]$ws.Cells.Replace($old, $labelChanges[$old])
    }
}

# ---------------------------------------------------------------
# 3. Updated figures after rerender
# ---------------------------------------------------------------

# rev_long sheet: FY2022 dollar figures
$revLong = $wb.Worksheets.Item("rev_long")
$revLong.Range("D544").Value = 1886.18847954
$revLong.Range("D551").Value = 2700.23621584

# Table 2 sheet: percent-change figures
$table2 = $wb.Worksheets.Item("Table 2")
$table2.Range("C11").Value = -4.32
$table2.Range("D11").Value = 7.89

# Revenue CAGR sheet (formerly Table 4.a)
$revenueCagr = $wb.Worksheets.Item("Revenue CAGR")
$revenueCagr.Range("F2").Value = 6.29
$revenueCagr.Range("B13").Value = -4.32
$revenueCagr.Range("C13").Value = 15.28
$revenueCagr.Range("D13").Value = 16.98
$revenueCagr.Range("E13").Value = 9.34
$revenueCagr.Range("F13").Value = 6.27
$revenueCagr.Range("G13").Value = 7.89

# year_totals sheet: FY2022 revenue and fiscal gap totals
$yearTotals = $wb.Worksheets.Item("year_totals")
$yearTotals.Range("C2").Value = 104544.13135584
$yearTotals.Range("D2").Value = 2715
